$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.593.68'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.812.00'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''226.47'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '''0.599'
$ws.Range('E6').Value = '  +3.24%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''37.79'
$ws.Range('E8').Value = '  +8.25%  '
$ws.Range('E9').Value = '  -3.51%  '
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('D11').Value = '''0.0972'
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('D12').Value = '2.072.89'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '''11.34'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = '1.823.63'
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').Value = '34.544.83'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').Value = '''68.82'
$ws.Range('D19').Value = '''244.67'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').Value = '''11.25'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '''4.15'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').Value = '''2.22'
$ws.Range('E24').Value = '  +4.89%  '
$ws.Range('D25').Value = '''172.25'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '''7.92'
$ws.Range('E26').Value = '  -1.00%  '
$ws.Range('D27').Value = '''17.37'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '''3.94'
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').Value = '''1.24'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').Value = '''0.0523'
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('D35').Value = '1.365.88'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = '''0.658'
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '''2.35'
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  +8.57%  '
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('D42').Value = '''81.16'
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '''0.942'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '''2.78'
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('D45').Value = '''14.07'
$ws.Range('E45').Value = '  +5.48%  '
$ws.Range('E46').Value = '  -2.47%  '
$ws.Range('D47').Value = '1.973.45'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''5.82'
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '''102.85'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = '0.0₆0121'
$ws.Range('E51').Value = '  -7.52%  '

Write-Output "Applied 84 cell updates to cryptos sheet"
